$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 corresponds to ZDBDID 254938-4 "Geologija rudnych mestoroždenij"
# Commit: "fix check for current holdings"
#  - Unikal (W5): was "x" -> now empty (no longer unique / currently held elsewhere too)
#  - Anzahl-FL-Bibliotheken (X5): was "0" -> now "1"
#  - FL-Bibliotheken (Y5): was "" -> now "DE-104"
$ws.Range("W5").Value = ""
$ws.Range("X5").Value = "'1"
$ws.Range("Y5").Value = "DE-104"
